$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$old,
        [string]$new
    )
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "FAILED to find:" $old
    } else {
        Write-Host "OK replaced:" $old
    }
    return $ok
}

# 1. "description, name " + "etc" (spell-checked as two runs) -> a single run,
#    same visible text "description, name etc".
Replace-Text "description, name etc" "description, name etc"

# 2. " trip (" + "eg" + ": a city, a restaurant)" (spell-checked) -> single run,
#    same visible text " trip (eg: a city, a restaurant)".
Replace-Text "trip (eg: a city, a restaurant)" "trip (eg: a city, a restaurant)"

# 3. "...booking a hotel room " + "etc" + ")..." (spell-checked) -> single run,
#    same visible text "...booking a hotel room etc). One of the most important ".
Replace-Text "booking a hotel room etc). One of the most important" "booking a hotel room etc). One of the most important"

# 4. "to your trip" -> "to their trip"
Replace-Text "to your trip" "to their trip"

# 5. "Once you've decided to share the trip with the people you're going, they can mak"
#    -> "Once someone has decided to share the trip with the people they're going, all of them can mak"
Replace-Text "Once you’ve decided to share the trip with the people you’re going, they can mak" "Once someone has decided to share the trip with the people they’re going, all of them can mak"

# 6. "Email notifications when someone added you to a trip"
#    -> "Send email notifications once someone's been added to a trip"
Replace-Text "Email notifications when someone added you to a trip" "Send email notifications once someone’s been added to a trip"

# 7. "Here you can make ... helps you make sure you will get to see everything ..."
#    -> "The user can make ... attraction they want ... helps people make sure they will get ..."
Replace-Text "Here you can make" "The user can make"
Replace-Text "touristic attraction you want to visit" "touristic attraction they want to visit"
Replace-Text "Putting them on a calendar helps you make sure you will get to see everything" "Putting them on a calendar helps people make sure they will get to see everything"

# 8. "A forum " + "section" (grammar-checked) -> single run, same visible text "A forum section".
Replace-Text "A forum section" "A forum section"

Write-Host "All replacements attempted."
